# Auto-generated script applying scheduled market-data updates to Gilgamesh_Profits workbook.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) per sheet/row,
# matching the latest Universalis market snapshot used by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4006
$ws.Range("I62").Value = 4013.5
$ws.Range("K62").Value = 4013.5
$ws.Range("M62").Value = -3389.5
$ws.Range("H65").Value = 4006
$ws.Range("I65").Value = 4013.5
$ws.Range("K65").Value = 20067.5
$ws.Range("M65").Value = -16947.5
$ws.Range("H70").Value = 5081.5
$ws.Range("J70").Value = 4934.6
$ws.Range("L70").Value = 14803.8
$ws.Range("N70").Value = -15343.8
$ws.Range("H73").Value = 5081.5
$ws.Range("J73").Value = 4934.6
$ws.Range("L73").Value = 14803.8
$ws.Range("N73").Value = -16675.8
$ws.Range("H107").Value = 401.33334
$ws.Range("I107").Value = 401.33334
$ws.Range("K107").Value = 401.33334
$ws.Range("M107").Value = 1518.66666
$ws.Range("H132").Value = 3181.5386
$ws.Range("I132").Value = 3572.111
$ws.Range("K132").Value = 10716.333
$ws.Range("M132").Value = -8186.332999999999
$ws.Range("H137").Value = 1670822.6
$ws.Range("I137").Value = 2177374.5
$ws.Range("J137").Value = 6437.5713
$ws.Range("K137").Value = 6532123.5
$ws.Range("L137").Value = 19312.7139
$ws.Range("M137").Value = -6529573.5
$ws.Range("N137").Value = -24412.7139
$ws.Range("H138").Value = 2045.9038
$ws.Range("I138").Value = 687.4
$ws.Range("J138").Value = 3303.7778
$ws.Range("K138").Value = 2062.2
$ws.Range("L138").Value = 9911.3334
$ws.Range("M138").Value = 3077.8
$ws.Range("N138").Value = -20191.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 31802.357
$ws.Range("I45").Value = 34061.77
$ws.Range("K45").Value = 34061.77
$ws.Range("M45").Value = -33684.77
$ws.Range("H74").Value = 266451.75
$ws.Range("I74").Value = 348811.7
$ws.Range("K74").Value = 348811.7
$ws.Range("M74").Value = -347937.7
$ws.Range("H77").Value = 266451.75
$ws.Range("I77").Value = 348811.7
$ws.Range("K77").Value = 1744058.5
$ws.Range("M77").Value = -1739690.5
$ws.Range("H88").Value = 2075.7
$ws.Range("I88").Value = 798.1
$ws.Range("J88").Value = 2714.5
$ws.Range("K88").Value = 798.1
$ws.Range("L88").Value = 2714.5
$ws.Range("M88").Value = -392.1
$ws.Range("N88").Value = -3526.5
$ws.Range("H91").Value = 2075.7
$ws.Range("I91").Value = 798.1
$ws.Range("J91").Value = 2714.5
$ws.Range("K91").Value = 798.1
$ws.Range("L91").Value = 2714.5
$ws.Range("M91").Value = 605.9
$ws.Range("N91").Value = -5522.5
$ws.Range("H97").Value = 1160.3334
$ws.Range("I97").Value = 1051.96
$ws.Range("K97").Value = 1051.96
$ws.Range("M97").Value = -555.96
$ws.Range("H102").Value = 3303.2942
$ws.Range("I102").Value = 3043.7334
$ws.Range("K102").Value = 3043.7334
$ws.Range("M102").Value = -1421.7334
$ws.Range("H110").Value = 758.8823
$ws.Range("I110").Value = 734.53845
$ws.Range("K110").Value = 734.53845
$ws.Range("M110").Value = 1310.46155
$ws.Range("H132").Value = 2188.568
$ws.Range("I132").Value = 1907.9395
$ws.Range("K132").Value = 5723.818499999999
$ws.Range("M132").Value = -3193.818499999999
$ws.Range("H139").Value = 78583.8
$ws.Range("J139").Value = 78583.8
$ws.Range("L139").Value = 78583.8
$ws.Range("N139").Value = -88863.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 32262926
$ws.Range("I94").Value = 43483890
$ws.Range("J94").Value = 2657.375
$ws.Range("K94").Value = 43483890
$ws.Range("L94").Value = 2657.375
$ws.Range("M94").Value = -43483439
$ws.Range("N94").Value = -3559.375
$ws.Range("H134").Value = 2511.762
$ws.Range("I134").Value = 1961.7812
$ws.Range("J134").Value = 4271.7
$ws.Range("K134").Value = 5885.3436
$ws.Range("L134").Value = 12815.1
$ws.Range("M134").Value = -3350.3436
$ws.Range("N134").Value = -17885.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 4729.5
$ws.Range("I22").Value = 5056.5
$ws.Range("J22").Value = 3748.5
$ws.Range("K22").Value = 5056.5
$ws.Range("L22").Value = 3748.5
$ws.Range("M22").Value = -4706.5
$ws.Range("N22").Value = -4448.5
$ws.Range("H31").Value = 3209506.2
$ws.Range("I31").Value = 3603.5
$ws.Range("K31").Value = 3603.5
$ws.Range("M31").Value = -3308.5
$ws.Range("H34").Value = 3209506.2
$ws.Range("I34").Value = 3603.5
$ws.Range("K34").Value = 3603.5
$ws.Range("M34").Value = -3401.5
$ws.Range("H58").Value = 2971.2917
$ws.Range("I58").Value = 2507.3333
$ws.Range("J58").Value = 3435.25
$ws.Range("K58").Value = 2507.3333
$ws.Range("L58").Value = 3435.25
$ws.Range("M58").Value = -2304.3333
$ws.Range("N58").Value = -3841.25
$ws.Range("H86").Value = 7778.933
$ws.Range("I86").Value = 7745
$ws.Range("K86").Value = 7745
$ws.Range("M86").Value = -6622
$ws.Range("H89").Value = 7778.933
$ws.Range("I89").Value = 7745
$ws.Range("K89").Value = 38725
$ws.Range("M89").Value = -33109
$ws.Range("H107").Value = 4167434
$ws.Range("I107").Value = 6250389
$ws.Range("K107").Value = 6250389
$ws.Range("M107").Value = -6248469
$ws.Range("H132").Value = 18527396
$ws.Range("I132").Value = 8217.076999999999
$ws.Range("K132").Value = 24651.231
$ws.Range("M132").Value = -22121.231
$ws.Range("H136").Value = 2971.2917
$ws.Range("I136").Value = 2507.3333
$ws.Range("J136").Value = 3435.25
$ws.Range("K136").Value = 7521.999899999999
$ws.Range("L136").Value = 10305.75
$ws.Range("M136").Value = -4971.999899999999
$ws.Range("N136").Value = -15405.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1769.0344
$ws.Range("I139").Value = 1423.9048
$ws.Range("K139").Value = 4271.7144
$ws.Range("M139").Value = 868.2856000000002
$ws.Range("H140").Value = 7115.757
$ws.Range("I140").Value = 3547.1765
$ws.Range("K140").Value = 10641.5295
$ws.Range("M140").Value = -5461.529500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 6322.9287
$ws.Range("J107").Value = 7856.4546
$ws.Range("L107").Value = 7856.4546
$ws.Range("N107").Value = -11696.4546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 30025
$ws.Range("I42").Value = 30025
$ws.Range("K42").Value = 30025
$ws.Range("M42").Value = -29462
$ws.Range("H43").Value = 10526.315
$ws.Range("I43").Value = 10000
$ws.Range("K43").Value = 10000
$ws.Range("M43").Value = -9807
$ws.Range("H45").Value = 30041
$ws.Range("I45").Value = 30041
$ws.Range("K45").Value = 30041
$ws.Range("M45").Value = -29634
$ws.Range("H47").Value = 40000
$ws.Range("J47").Value = 40000
$ws.Range("L47").Value = 40000
$ws.Range("N47").Value = -40980
$ws.Range("H49").Value = 30025
$ws.Range("I49").Value = 30025
$ws.Range("K49").Value = 30025
$ws.Range("M49").Value = -29878
$ws.Range("H52").Value = 40000
$ws.Range("J52").Value = 40000
$ws.Range("L52").Value = 40000
$ws.Range("N52").Value = -40466
$ws.Range("H68").Value = 5018.1
$ws.Range("I68").Value = 3795.5
$ws.Range("J68").Value = 5833.1665
$ws.Range("K68").Value = 3795.5
$ws.Range("L68").Value = 5833.1665
$ws.Range("M68").Value = -3046.5
$ws.Range("N68").Value = -7331.1665
$ws.Range("H71").Value = 5018.1
$ws.Range("I71").Value = 3795.5
$ws.Range("J71").Value = 5833.1665
$ws.Range("K71").Value = 18977.5
$ws.Range("L71").Value = 29165.8325
$ws.Range("M71").Value = -15233.5
$ws.Range("N71").Value = -36653.8325
$ws.Range("H132").Value = 4261.9707
$ws.Range("I132").Value = 1750.2693
$ws.Range("K132").Value = 5250.8079
$ws.Range("M132").Value = -2720.8079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 30028
$ws.Range("I49").Value = 30056
$ws.Range("J49").Value = 30000
$ws.Range("K49").Value = 30056
$ws.Range("L49").Value = 30000
$ws.Range("M49").Value = -29826
$ws.Range("N49").Value = -30460
$ws.Range("H62").Value = 6158
$ws.Range("I62").Value = 5724
$ws.Range("K62").Value = 5724
$ws.Range("M62").Value = -5100
$ws.Range("H65").Value = 6158
$ws.Range("I65").Value = 5724
$ws.Range("K65").Value = 28620
$ws.Range("M65").Value = -25500
$ws.Range("H113").Value = 906.7
$ws.Range("I113").Value = 735.2727
$ws.Range("K113").Value = 2205.8181
$ws.Range("M113").Value = -35.81809999999996
$ws.Range("H132").Value = 2819.5806
$ws.Range("I132").Value = 2867.037
$ws.Range("K132").Value = 8601.110999999999
$ws.Range("M132").Value = -6071.110999999999
